# Fixed stacked bar plot to include sex, corrected truth input in stacked_plot()
# - Adds a new "SEX" worksheet (age_group x Female/Male breakdown), inserted
#   right after "Hypertension (1=yes; 0=no)" and before "All Combined".
# - Extends the "All Combined" sheet with the Female/Male/Total (+ %) columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "SEX" worksheet after "Hypertension (1=yes; 0=no)"
# ---------------------------------------------------------------------------
$hyp = $wb.Worksheets.Item("Hypertension (1=yes; 0=no)")
$sexSheet = $wb.Worksheets.Add($null, $hyp)
$sexSheet.Name = "SEX"

function Set-HeaderStyle($cell) {
    $cell.Borders.LineStyle = 1
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$sexHeaders = @("age_group", "Female", "Male", "Total", "Female_%", "Male_%", "Total_%")
for ($col = 1; $col -le $sexHeaders.Length; $col++) {
    $cell = $sexSheet.Cells.Item(1, $col)
    $cell.Value = $sexHeaders[$col - 1]
    Set-HeaderStyle $cell
}

$sexRows = @(
    @("18-29", 8, 8, 16, 50, 50, 100),
    @("30-39", 16, 12, 28, 56.99999999999999, 43, 100),
    @("40-49", 27, 51, 78, 35, 65, 100),
    @("50-59", 38, 66, 104, 37, 63, 100),
    @("60-69", 63, 109, 172, 37, 63, 100),
    @("70-79", 73, 134, 207, 35, 65, 100),
    @("80-89", 53, 81, 134, 40, 60, 100),
    @("90-99", 3, 1, 4, 75, 25, 100),
    @("Total", 281, 462, 743, 38, 62, 100)
)

$rowIndex = 2
foreach ($row in $sexRows) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $cell = $sexSheet.Cells.Item($rowIndex, $col)
        $cell.Value = $row[$col - 1]
        if ($col -eq 1) {
            Set-HeaderStyle $cell
        }
    }
    $rowIndex++
}

# ---------------------------------------------------------------------------
# 2. Extend "All Combined" with the Female/Male columns (T:Y)
# ---------------------------------------------------------------------------
$combined = $wb.Worksheets.Item("All Combined")

$combinedHeaders = @("Female", "Male", "Total", "Female_%", "Male_%", "Total_%")
for ($i = 0; $i -lt $combinedHeaders.Length; $i++) {
    $cell = $combined.Cells.Item(1, 20 + $i)
    $cell.Value = $combinedHeaders[$i]
    Set-HeaderStyle $cell
}

$combinedRows = @(
    @(8, 8, 16, 50, 50, 100),
    @(16, 12, 28, 56.99999999999999, 43, 100),
    @(27, 51, 78, 35, 65, 100),
    @(38, 66, 104, 37, 63, 100),
    @(63, 109, 172, 37, 63, 100),
    @(73, 134, 207, 35, 65, 100),
    @(53, 81, 134, 40, 60, 100),
    @(3, 1, 4, 75, 25, 100),
    @(281, 462, 743, 38, 62, 100)
)

$rowIndex = 2
foreach ($row in $combinedRows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $combined.Cells.Item($rowIndex, 20 + $i).Value = $row[$i]
    }
    $rowIndex++
}
